# Add season record columns (Wins, Losses, Ties) to the right of the
# existing statistics table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - mirror the formatting already used by the other
# header cells (bold font, thin border around, centered alignment) by
# copying the format from an existing header cell.
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-54: every team row gets the same season record.
$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("AD$r").Value = 64
    $ws.Range("AE$r").Value = 98
    $ws.Range("AF$r").Value = 0
}

Write-Output "done"
